# Updates the "cryptos" price-tracker sheet in place: refreshed Price (D)
# and Volume(1h) (E) figures for every coin, plus three pairs/triples of
# rows (17/18, 25/26, 34/35/36) whose Coin name + Link swapped position
# because the underlying ranking reordered.
#
# Each entry below only carries the columns that actually changed for
# that row (column A's row index and unaffected B/C/D/E values are left
# completely alone).

$rows = @(
    @{ Row = 2;  D = "67.522.15";  E = "  -0.93%  " },
    @{ Row = 3;  D = "3.640.27";   E = "  -1.32%  " },
    @{ Row = 4;  D = "1.01";       E = "  +0.54%  " },
    @{ Row = 5;  D = "586.47";     E = "  -1.39%  " },
    @{ Row = 6;  D = "185.22";     E = "  +2.46%  " },
    @{ Row = 7;  E = "  -2.16%  " },
    @{ Row = 8;  E = "  -0.29%  " },
    @{ Row = 9;  D = "0.684";      E = "  -4.18%  " },
    @{ Row = 10; E = "  -8.42%  " },
    @{ Row = 11; D = "54.90";      E = "  -2.26%  " },
    @{ Row = 12; E = "  -10.54%  " },
    @{ Row = 13; D = "10.05";      E = "  -2.99%  " },
    @{ Row = 14; D = "4.232.55";   E = "  -0.89%  " },
    @{ Row = 15; D = "3.649.35";   E = "  -0.90%  " },
    @{ Row = 17; B = "Chainlink";                 C = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link";          D = "18.56";     E = "  -3.75%  " },
    @{ Row = 18; B = "WrappedBTC";                C = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc";          D = "67.419.51"; E = "  -0.85%  " },
    @{ Row = 19; D = "1.09";       E = "  -2.93%  " },
    @{ Row = 20; D = "12.36";      E = "  -3.33%  " },
    @{ Row = 21; D = "398.47";     E = "  -2.41%  " },
    @{ Row = 22; D = "4.33";       E = "  -4.99%  " },
    @{ Row = 23; D = "86.43" },
    @{ Row = 24; E = "  -4.12%  " },
    @{ Row = 25; B = "InternetComputer(DFINITY)";  C = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D = "12.33";     E = "  -3.51%  " },
    @{ Row = 26; B = "RenderToken";                C = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr";       D = "10.63";     E = "  -2.92%  " },
    @{ Row = 27; D = "6.06";       E = "  +0.12%  " },
    @{ Row = 28; D = "3.63";       E = "  -5.96%  " },
    @{ Row = 29; D = "9.14";       E = "  -2.86%  " },
    @{ Row = 30; D = "31.51";      E = "  -3.77%  " },
    @{ Row = 31; D = "6.93";       E = "  -4.09%  " },
    @{ Row = 32; D = "67.30";      E = "  +4.64%  " },
    @{ Row = 33; D = "12.01";      E = "  -3.32%  " },
    @{ Row = 34; B = "InjectiveProtocol"; C = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; D = "43.04";  E = "  -0.74%  " },
    @{ Row = 35; B = "Hedera";            C = "https://coinranking.com/coin/jad286TjB+hedera-hbar";           D = "0.113"; E = "  -3.01%  " },
    @{ Row = 36; B = "Bittensor";         C = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao";         D = "591.01"; E = "  -2.09%  " },
    @{ Row = 37; E = "  -0.07%  " },
    @{ Row = 38; E = "  +0.14%  " },
    @{ Row = 39; E = "  -4.15%  " },
    @{ Row = 40; E = "  -0.78%  " },
    @{ Row = 41; D = "0.0₃0733";   E = "  -17.17%  " },
    @{ Row = 42; D = "2.83";       E = "  -5.44%  " },
    @{ Row = 43; D = "0.0415";     E = "  -5.04%  " },
    @{ Row = 44; E = "  +1.52%  " },
    @{ Row = 45; D = "2.45";       E = "  -12.49%  " },
    @{ Row = 46; E = "  -0.64%  " },
    @{ Row = 47; D = "2.717.24";   E = "  -0.69%  " },
    @{ Row = 48; D = "141.24";     E = "  -0.72%  " },
    @{ Row = 49; D = "8.56";       E = "  -6.74%  " },
    @{ Row = 50; E = "  -5.56%  " },
    @{ Row = 51; E = "  -4.68%  " }
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($entry in $rows) {
    $r = $entry.Row

    if ($entry.ContainsKey("B")) {
        $ws.Range("B$r").Value = $entry.B
    }
    if ($entry.ContainsKey("C")) {
        $ws.Range("C$r").Value = $entry.C
    }
    if ($entry.ContainsKey("D")) {
        # Many Price strings look numeric ("1.01", "586.47", ...). A plain
        # Value assignment would let Excel silently coerce them to a
        # Number, which is not what the sheet stores (these stay text,
        # e.g. "67.522.15" / "0.0₃0733" would be mangled otherwise too).
        # Force Text format, assign, then clear the temporary formatting
        # so the cell keeps its original (default) style.
        $cell = $ws.Range("D$r")
        $cell.NumberFormat = "@"
        $cell.Value = $entry.D
        $cell.ClearFormats()
    }
    if ($entry.ContainsKey("E")) {
        $ws.Range("E$r").Value = $entry.E
    }
}
